$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 16-18 (table shrinks from 18 rows to 15 rows)
$ws.Range("A16:C18").Delete()

$data = @(
    @("<zero>", "<zero>", 17),
    @("<part>", "<part>", 19),
    @("<zulu>", "<zulu>", 15),
    @("<water>", "<water>", 23),
    @("<can>", "<can>", 15),
    @("<a>", "<a>", 22),
    @("<you>", "<you>", 19),
    @("<number>", "<nomeo>", 9),
    @("<four>", "<for>", 22),
    @("<word>", "<word>", 16),
    @("<tango>", "<tango>", 24),
    @("<backspace>", "<backspace>", 21),
    @("<said>", "<said>", 25),
    @("<november>", "<tangr>", 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
